$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swap: Litecoin (row 26) <-> NEARProtocol (row 27) ---
# New values after the swap (rank numbers in column A stay put):
# Row 26 becomes NEARProtocol, Row 27 becomes Litecoin
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'6.15"
$ws.Range("E26").Value = "  +7.99%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'98.53"
$ws.Range("E27").Value = "  +10.69%  "

# --- Row swap: Dai (row 32) <-> InternetComputer(DFINITY) (row 33) ---
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'11.33"
$ws.Range("E32").Value = "  +6.59%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.23%  "

# --- Simple value updates (Price / Volume columns) ---
# Note: Price values that look like plain numbers (single decimal point) are
# entered with a leading apostrophe so Excel keeps them as literal text
# (matching the source data, which stores prices as text, not numbers).
# Prices already containing two dots (e.g. "97.757.71") are naturally kept
# as text by Excel without needing the apostrophe.
$ws.Range("D2").Value = "97.757.71"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.420.41"
$ws.Range("E3").Value = "  +3.95%  "
$ws.Range("D5").Value = "'254.97"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'652.86"
$ws.Range("E6").Value = "  +4.35%  "
$ws.Range("D7").Value = "'1.48"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "'0.427"
$ws.Range("E8").Value = "  +6.42%  "
$ws.Range("E9").Value = "  +8.57%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.418.53"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("D12").Value = "'0.212"
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("D13").Value = "'41.64"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "'6.30"
$ws.Range("E14").Value = "  +14.97%  "
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "97.471.15"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "4.058.19"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("E18").Value = "  +34.05%  "
$ws.Range("D19").Value = "3.424.31"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("D20").Value = "'17.53"
$ws.Range("E20").Value = "  +12.38%  "
$ws.Range("D21").Value = "'0.494"
$ws.Range("E21").Value = "  +43.12%  "
$ws.Range("D22").Value = "'10.71"
$ws.Range("E22").Value = "  +14.14%  "
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'503.08"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D28").Value = "'12.60"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").Value = "3.600.59"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("D30").Value = "'0.152"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "'0.201"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'0.574"
$ws.Range("E35").Value = "  +19.58%  "
$ws.Range("D36").Value = "'29.74"
$ws.Range("E36").Value = "  +6.40%  "
$ws.Range("D37").Value = "'2.25"
$ws.Range("E37").Value = "  +15.14%  "
$ws.Range("E38").Value = "  +5.80%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "'1.43"
$ws.Range("E40").Value = "  +14.75%  "
$ws.Range("D41").Value = "'518.13"
$ws.Range("E41").Value = "  +5.00%  "
$ws.Range("D42").Value = "'24.73"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'0.873"
$ws.Range("E43").Value = "  +12.00%  "
$ws.Range("D44").Value = "'3.72"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").Value = "'0.0418"
$ws.Range("E45").Value = "  +22.91%  "
$ws.Range("D46").Value = "'5.50"
$ws.Range("E46").Value = "  +14.11%  "
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").Value = "'8.24"
$ws.Range("E48").Value = "  +12.58%  "
$ws.Range("E50").Value = "  +12.20%  "
$ws.Range("D51").Value = "'2.05"
$ws.Range("E51").Value = "  +5.19%  "
